$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "1.000", "0.07690")
# keep their exact formatting instead of being coerced into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.576.39'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").Value = '1.918.81'
$ws.Range("E3").Value = '  -0.99%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").Value = '246.06'
$ws.Range("E5").Value = '  +2.31%  '
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("D7").Value = '0.4735'
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").Value = '0.2883'
$ws.Range("E8").Value = '  +0.34%  '
$ws.Range("D9").Value = '0.06825'
$ws.Range("E9").Value = '  +3.52%  '
$ws.Range("D10").Value = '104.82'
$ws.Range("E10").Value = '  -3.32%  '
$ws.Range("D11").Value = '18.29'
$ws.Range("E11").Value = '  -6.50%  '
$ws.Range("D12").Value = '1.920.83'
$ws.Range("E12").Value = '  -0.72%  '
$ws.Range("D13").Value = '0.07690'
$ws.Range("E13").Value = '  +1.16%  '
$ws.Range("D14").Value = '5.257'
$ws.Range("E14").Value = '  +2.24%  '
$ws.Range("D15").Value = '0.6682'
$ws.Range("E15").Value = '  +2.05%  '
$ws.Range("D16").Value = '291.32'
$ws.Range("E16").Value = '  -8.53%  '
$ws.Range("D17").Value = '30.573.09'
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D18").Value = '0.9994'
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.000007580'
$ws.Range("E19").Value = '  +0.51%  '
$ws.Range("D20").Value = '12.91'
$ws.Range("E20").Value = '  -1.12%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.171.16'
$ws.Range("E21").Value = '  +1.13%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '5.507'
$ws.Range("E22").Value = '  +5.59%  '
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("D24").Value = '6.347'
$ws.Range("E24").Value = '  -0.99%  '
$ws.Range("D25").Value = '9.382'
$ws.Range("E25").Value = '  +0.38%  '
$ws.Range("D26").Value = '167.83'
$ws.Range("E26").Value = '  +1.33%  '
$ws.Range("D27").Value = '21.01'
$ws.Range("E27").Value = '  +2.49%  '
$ws.Range("D28").Value = '2.110'
$ws.Range("E28").Value = '  +2.70%  '
$ws.Range("D29").Value = '0.1061'
$ws.Range("E29").Value = '  -5.14%  '
$ws.Range("D30").Value = '1.396'
$ws.Range("E30").Value = '  +3.70%  '
$ws.Range("D31").Value = '4.158'
$ws.Range("E31").Value = '  +0.84%  '
$ws.Range("D32").Value = '4.056'
$ws.Range("E32").Value = '  +2.10%  '
$ws.Range("D33").Value = '0.05027'
$ws.Range("E33").Value = '  -0.85%  '
$ws.Range("D34").Value = '0.7335'
$ws.Range("E34").Value = '  -1.63%  '
$ws.Range("D35").Value = '1.140'
$ws.Range("E35").Value = '  -1.91%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '2.744'
$ws.Range("E36").Value = '  +1.12%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.02051'
$ws.Range("E37").Value = '  +3.73%  '
$ws.Range("D38").Value = '2.686'
$ws.Range("E38").Value = '  -1.10%  '
$ws.Range("D39").Value = '2.046'
$ws.Range("E39").Value = '  +0.78%  '
$ws.Range("D40").Value = '110.74'
$ws.Range("E40").Value = '  +2.87%  '
$ws.Range("D41").Value = '0.8750'
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("D42").Value = '0.4380'
$ws.Range("E42").Value = '  +4.90%  '
$ws.Range("D43").Value = '5.856'
$ws.Range("E43").Value = '  -0.85%  '
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("D45").Value = '66.98'
$ws.Range("E45").Value = '  -2.20%  '
$ws.Range("D46").Value = '7.237'
$ws.Range("E46").Value = '  -1.07%  '
$ws.Range("D47").Value = '9.295'
$ws.Range("E47").Value = '  -0.39%  '
$ws.Range("D48").Value = '47.91'
$ws.Range("E48").Value = '  +8.82%  '
$ws.Range("D49").Value = '0.1226'
$ws.Range("E49").Value = '  +1.08%  '
$ws.Range("D50").Value = '34.79'
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").Value = '0.4019'
$ws.Range("E51").Value = '  +4.14%  '
